# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume 1h) updates for rows 2-49.
# Price/volume values are textual (e.g. "28.315.51", "  +4.03%  ") so the
# target cells must stay text rather than being auto-coerced to numbers.
$updates = @{
    2  = @{ D = "28.315.51"; E = "  +4.03%  " }
    3  = @{ D = "1.732.37";  E = "  +2.70%  " }
    4  = @{ E = "  -0.06%  " }
    5  = @{ D = "219.56" }
    6  = @{ D = "0.524";    E = "  +0.70%  " }
    7  = @{ E = "  -0.12%  " }
    8  = @{ D = "24.16";    E = "  +6.28%  " }
    9  = @{ E = "  +3.14%  " }
    10 = @{ D = "0.0637";   E = "  +1.64%  " }
    11 = @{ D = "0.0895";   E = "  +0.27%  " }
    12 = @{ D = "1.977.89"; E = "  +2.69%  " }
    13 = @{ D = "1.722.84"; E = "  +1.99%  " }
    14 = @{ E = "  +2.11%  " }
    15 = @{ D = "0.565";    E = "  +1.99%  " }
    16 = @{ D = "67.73";    E = "  +0.82%  " }
    17 = @{ D = "28.316.92"; E = "  +4.00%  " }
    18 = @{ D = "242.88";   E = "  +1.54%  " }
    19 = @{ E = "  +1.26%  " }
    20 = @{ D = "7.96";     E = "  -2.60%  " }
    21 = @{ E = "  -0.03%  " }
    22 = @{ D = "4.67";     E = "  +2.23%  " }
    23 = @{ D = "9.79";     E = "  +1.76%  " }
    24 = @{ D = "2.11";     E = "  -0.20%  " }
    25 = @{ E = "  +0.86%  " }
    26 = @{ D = "7.55";     E = "  +3.39%  " }
    27 = @{ D = "16.65";    E = "  +0.90%  " }
    28 = @{ E = "  +0.73%  " }
    29 = @{ E = "  -0.27%  " }
    30 = @{ D = "0.0518";   E = "  +3.21%  " }
    31 = @{ E = "  +2.58%  " }
    32 = @{ E = "  +0.65%  " }
    33 = @{ D = "3.29";     E = "  +1.26%  " }
    34 = @{ D = "1.489.18"; E = "  -5.46%  " }
    35 = @{ E = "  -1.72%  " }
    36 = @{ D = "0.980";    E = "  +2.98%  " }
    37 = @{ D = "0.605";    E = "  +0.41%  " }
    38 = @{ E = "  +0.58%  " }
    39 = @{ E = "  +1.22%  " }
    40 = @{ D = "1.07";     E = "  +0.35%  " }
    41 = @{ D = "70.35";    E = "  +1.03%  " }
    42 = @{ E = "  -0.09%  " }
    43 = @{ E = "  +0.25%  " }
    44 = @{ D = "1.881.25"; E = "  +2.22%  " }
    45 = @{ E = "  +1.49%  " }
    46 = @{ D = "0.798";    E = "  +1.46%  " }
    47 = @{ E = "  +7.95%  " }
    48 = @{ E = "  +7.39%  " }
    49 = @{ D = "90.86";    E = "  -0.42%  " }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    if ($rowData.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $rowData["D"]
    }
    if ($rowData.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $rowData["E"]
    }
}

# Rows 50 and 51 swap their coin identity (Algorand <-> EnergySwap) along with
# updated price/volume values.
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.22"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.16%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.105"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.31%  "
